# Generate Report for Handoff
# Row 3 ("b.md") moves from "Handed back: in sync with en-US" to
# "Ready for handoff" with a fresh Latest Handoff File / Datetime, for
# both the zh-cn and de-de locales (and the Overview roll-up sheet).

$wb = $excel.ActiveWorkbook

# ---- Overview sheet --------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")
$ov.Range("B3").Value = "Ready for handoff"
$ov.Range("C3").Value = "Ready for handoff"
$ov.Range("D3").Value = "2016-32-20 14:32:41"

# ---- zh-cn sheet -------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")
$zh.Range("C3").Value = "Ready for handoff"
$zh.Range("D3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$zh.Range("E3").Value = "2016-03-20 14:32:38"
foreach ($h in $zh.Hyperlinks) {
    if ($h.Range.Address() -eq '$D$3') {
        $h.TextToDisplay = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
    }
}

# ---- de-de sheet -------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")
$de.Range("C3").Value = "Ready for handoff"
$de.Range("D3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$de.Range("E3").Value = "2016-03-20 14:32:41"
foreach ($h in $de.Hyperlinks) {
    if ($h.Range.Address() -eq '$D$3') {
        $h.TextToDisplay = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
    }
}
